# Evaluations.xlsx edit script
# Adds a new "MERGED RAG" evaluation block (gemini / llama / deepseek) and
# highlights the best scores in the existing "Quran RAG" table with
# bold / underline formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Highlight best-scoring cells in the existing table (rows 4-10)
# ---------------------------------------------------------------------

# Column E (Average BERTScore Recall) best score is on row 7 (gemini)
$ws.Range("E7:G7").Font.Bold = $true

# Row 9 (Mistral-7b temp=0.1) is best for B,H,K,N,Q - but NOT E (E stays as-is)
$ws.Range("B9:D9").Font.Bold = $true
$ws.Range("H9:S9").Font.Bold = $true

# Row 8 (llama-4-maverick) gets underlined instead
$ws.Range("B8:S8").Font.Underline = $true

# ---------------------------------------------------------------------
# 2. Normalize a couple of cells that carried stray formatting
# ---------------------------------------------------------------------
$ws.Range("A13").Font.Bold = $false
$ws.Range("T2").Font.Bold = $false
$ws.Range("W2").Font.Bold = $false

# ---------------------------------------------------------------------
# 3. Extend + bold the "Quran RAG" header row (row 12) across the
#    whole table width (A:AE) to match the newly added block below
# ---------------------------------------------------------------------
$ws.Range("X12:AE12").HorizontalAlignment = -4108
$ws.Range("B12:AE12").Font.Bold = $true
$ws.Range("A12:AE12").Merge()

# ---------------------------------------------------------------------
# 4. New "MERGED RAG" section
# ---------------------------------------------------------------------

# Header row 15
$ws.Range("A15").Value = "MERGED RAG"
$ws.Range("A15:AE15").HorizontalAlignment = -4108
$ws.Range("A15:AE15").Font.Bold = $true
$ws.Range("A15:AE15").Merge()

# Row 16 - gemini
$ws.Range("A16").Value = "gemini"
$ws.Range("B16").Value = 0.82879999999999998
$ws.Range("E16").Value = 0.80769999999999997
$ws.Range("H16").Value = 0.81789999999999996
$ws.Range("K16").Value = [double]"4.2099999999999999E-2"
$ws.Range("N16").Value = [double]"2.5000000000000001E-3"
$ws.Range("Q16").Value = [double]"3.85E-2"
$ws.Range("B16:S16").HorizontalAlignment = -4108

# Row 17 - llama
$ws.Range("A17").Value = "llama"
$ws.Range("B17").Value = 0.80840000000000001
$ws.Range("E17").Value = 0.80979999999999996
$ws.Range("H17").Value = 0.80879999999999996
$ws.Range("K17").Value = [double]"8.0199999999999994E-2"
$ws.Range("N17").Value = [double]"8.8999999999999999E-3"
$ws.Range("Q17").Value = [double]"5.3800000000000001E-2"
$ws.Range("B17:S17").HorizontalAlignment = -4108

# Row 18 - deepseek_r1_distill_llama_70b
$ws.Range("A18").Value = "deepseek_r1_distill_llama_70b"
$ws.Range("B18").Value = 0.79059999999999997
$ws.Range("E18").Value = 0.81940000000000002
$ws.Range("H18").Value = 0.80449999999999999
$ws.Range("K18").Value = 0.11550000000000001
$ws.Range("N18").Value = [double]"1.12E-2"
$ws.Range("Q18").Value = [double]"6.4100000000000004E-2"
$ws.Range("B18:S18").HorizontalAlignment = -4108

# Merge the 3-column groups for the 3 new data rows, same pattern as
# the rest of the table (B:D, E:G, H:J, K:M, N:P, Q:S)
$ws.Range("B16:D16").Merge()
$ws.Range("E16:G16").Merge()
$ws.Range("H16:J16").Merge()
$ws.Range("K16:M16").Merge()
$ws.Range("N16:P16").Merge()
$ws.Range("Q16:S16").Merge()

$ws.Range("B17:D17").Merge()
$ws.Range("E17:G17").Merge()
$ws.Range("H17:J17").Merge()
$ws.Range("K17:M17").Merge()
$ws.Range("N17:P17").Merge()
$ws.Range("Q17:S17").Merge()

$ws.Range("B18:D18").Merge()
$ws.Range("E18:G18").Merge()
$ws.Range("H18:J18").Merge()
$ws.Range("K18:M18").Merge()
$ws.Range("N18:P18").Merge()
$ws.Range("Q18:S18").Merge()

# ---------------------------------------------------------------------
# 5. Window / view state
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 108
$ws.Range("E18:G18").Select()
